$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.811.62"
$ws.Range("E2").Value = "  +8.43%  "
$ws.Range("D3").Value = "1.954.32"
$ws.Range("E3").Value = "  +6.82%  "
$ws.Range("E4").Value = "  -0.52%  "
$ws.Range("D5").Formula = '=TEXT("342.80","@")'
$ws.Range("D5").Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("E5").Value = "  +2.82%  "
$ws.Range("E6").Value = "  -0.31%  "
$ws.Range("D7").Formula = '=TEXT("0.4772","@")'
$ws.Range("D7").Copy()
$ws.Range("D7").PasteSpecial(-4163)
$ws.Range("E7").Value = "  +4.04%  "
$ws.Range("D8").Formula = '=TEXT("0.4151","@")'
$ws.Range("D8").Copy()
$ws.Range("D8").PasteSpecial(-4163)
$ws.Range("E8").Value = "  +8.32%  "
$ws.Range("D9").Formula = '=TEXT("48.10","@")'
$ws.Range("D9").Copy()
$ws.Range("D9").PasteSpecial(-4163)
$ws.Range("E9").Value = "  +3.95%  "
$ws.Range("D10").Formula = '=TEXT("0.08263","@")'
$ws.Range("D10").Copy()
$ws.Range("D10").PasteSpecial(-4163)
$ws.Range("E10").Value = "  +5.19%  "
$ws.Range("D11").Formula = '=TEXT("1.038","@")'
$ws.Range("D11").Copy()
$ws.Range("D11").PasteSpecial(-4163)
$ws.Range("E11").Value = "  +8.22%  "
$ws.Range("D12").Formula = '=TEXT("22.77","@")'
$ws.Range("D12").Copy()
$ws.Range("D12").PasteSpecial(-4163)
$ws.Range("E12").Value = "  +8.07%  "
$ws.Range("D13").Value = "1.954.69"
$ws.Range("E13").Value = "  +6.77%  "
$ws.Range("D14").Formula = '=TEXT("6.191","@")'
$ws.Range("D14").Copy()
$ws.Range("D14").PasteSpecial(-4163)
$ws.Range("E14").Value = "  +5.93%  "
$ws.Range("D15").Formula = '=TEXT("7.415","@")'
$ws.Range("D15").Copy()
$ws.Range("D15").PasteSpecial(-4163)
$ws.Range("E15").Value = "  +5.03%  "
$ws.Range("D16").Formula = '=TEXT("92.24","@")'
$ws.Range("D16").Copy()
$ws.Range("D16").PasteSpecial(-4163)
$ws.Range("E16").Value = "  +2.79%  "
$ws.Range("E17").Value = "  -0.52%  "
$ws.Range("E18").Value = "  +3.79%  "
$ws.Range("D19").Formula = '=TEXT("0.06695","@")'
$ws.Range("D19").Copy()
$ws.Range("D19").PasteSpecial(-4163)
$ws.Range("E19").Value = "  +1.62%  "
$ws.Range("D20").Formula = '=TEXT("18.05","@")'
$ws.Range("D20").Copy()
$ws.Range("D20").PasteSpecial(-4163)
$ws.Range("D21").Formula = '=TEXT("1.000","@")'
$ws.Range("D21").Copy()
$ws.Range("D21").PasteSpecial(-4163)
$ws.Range("E21").Value = "  -0.38%  "
$ws.Range("D22").Value = "29.774.97"
$ws.Range("E22").Value = "  +8.36%  "
$ws.Range("D23").Formula = '=TEXT("5.587","@")'
$ws.Range("D23").Copy()
$ws.Range("D23").PasteSpecial(-4163)
$ws.Range("E23").Value = "  +5.42%  "
$ws.Range("E24").Value = "  +4.45%  "
$ws.Range("D25").Formula = '=TEXT("2.265","@")'
$ws.Range("D25").Copy()
$ws.Range("D25").PasteSpecial(-4163)
$ws.Range("E25").Value = "  -0.32%  "
$ws.Range("D26").Value = "2.177.95"
$ws.Range("E26").Value = "  +6.16%  "
$ws.Range("D27").Formula = '=TEXT("162.03","@")'
$ws.Range("D27").Copy()
$ws.Range("D27").PasteSpecial(-4163)
$ws.Range("E27").Value = "  +1.85%  "
$ws.Range("D28").Formula = '=TEXT("20.20","@")'
$ws.Range("D28").Copy()
$ws.Range("D28").PasteSpecial(-4163)
$ws.Range("E28").Value = "  +4.27%  "
$ws.Range("E29").Value = "  +6.99%  "
$ws.Range("E30").Value = "  +7.65%  "
$ws.Range("D31").Formula = '=TEXT("122.98","@")'
$ws.Range("D31").Copy()
$ws.Range("D31").PasteSpecial(-4163)
$ws.Range("E31").Value = "  +4.30%  "
$ws.Range("D32").Formula = '=TEXT("1.010","@")'
$ws.Range("D32").Copy()
$ws.Range("D32").PasteSpecial(-4163)
$ws.Range("E32").Value = "  +8.51%  "
$ws.Range("D33").Formula = '=TEXT("0.09645","@")'
$ws.Range("D33").Copy()
$ws.Range("D33").PasteSpecial(-4163)
$ws.Range("E33").Value = "  +2.74%  "
$ws.Range("D34").Formula = '=TEXT("1.482","@")'
$ws.Range("D34").Copy()
$ws.Range("D34").PasteSpecial(-4163)
$ws.Range("E34").Value = "  +12.67%  "
$ws.Range("D35").Formula = '=TEXT("3.687","@")'
$ws.Range("D35").Copy()
$ws.Range("D35").PasteSpecial(-4163)
$ws.Range("E35").Value = "  +3.08%  "
$ws.Range("E36").Value = "  +5.86%  "
$ws.Range("D37").Formula = '=TEXT("0.06308","@")'
$ws.Range("D37").Copy()
$ws.Range("D37").PasteSpecial(-4163)
$ws.Range("E37").Value = "  +6.23%  "
$ws.Range("D38").Formula = '=TEXT("0.02316","@")'
$ws.Range("D38").Copy()
$ws.Range("D38").PasteSpecial(-4163)
$ws.Range("E38").Value = "  +5.90%  "
$ws.Range("D39").Formula = '=TEXT("8.509","@")'
$ws.Range("D39").Copy()
$ws.Range("D39").PasteSpecial(-4163)
$ws.Range("E39").Value = "  +5.07%  "
$ws.Range("E40").Value = "  +3.80%  "
$ws.Range("D41").Formula = '=TEXT("0.6105","@")'
$ws.Range("D41").Copy()
$ws.Range("D41").PasteSpecial(-4163)
$ws.Range("E41").Value = "  +6.50%  "
$ws.Range("E42").Value = "  +8.61%  "
$ws.Range("B43").Value = "Algorand"
$ws.Range("C43").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D43").Formula = '=TEXT("0.1897","@")'
$ws.Range("D43").Copy()
$ws.Range("D43").PasteSpecial(-4163)
$ws.Range("E43").Value = "  +4.09%  "
$ws.Range("B44").Value = "Frax"
$ws.Range("C44").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D44").Formula = '=TEXT("1.000","@")'
$ws.Range("D44").Copy()
$ws.Range("D44").PasteSpecial(-4163)
$ws.Range("E44").Value = "  -0.34%  "
$ws.Range("D45").Formula = '=TEXT("2.397","@")'
$ws.Range("D45").Copy()
$ws.Range("D45").PasteSpecial(-4163)
$ws.Range("E45").Value = "  +34.21%  "
$ws.Range("D46").Formula = '=TEXT("1.272","@")'
$ws.Range("D46").Copy()
$ws.Range("D46").PasteSpecial(-4163)
$ws.Range("E46").Value = "  -0.89%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Formula = '=TEXT("12.57","@")'
$ws.Range("D47").Copy()
$ws.Range("D47").PasteSpecial(-4163)
$ws.Range("E47").Value = "  +6.24%  "
$ws.Range("B48").Value = "Decentraland"
$ws.Range("C48").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D48").Formula = '=TEXT("0.5727","@")'
$ws.Range("D48").Copy()
$ws.Range("D48").PasteSpecial(-4163)
$ws.Range("E48").Value = "  +6.29%  "
$ws.Range("D49").Formula = '=TEXT("1.987","@")'
$ws.Range("D49").Copy()
$ws.Range("D49").PasteSpecial(-4163)
$ws.Range("E49").Value = "  +5.26%  "
$ws.Range("E50").Value = "  +6.90%  "
$ws.Range("D51").Formula = '=TEXT("113.69","@")'
$ws.Range("D51").Copy()
$ws.Range("D51").PasteSpecial(-4163)
$ws.Range("E51").Value = "  +3.00%  "
$excel.CutCopyMode = $false
